# Release 1.6.3 update: rename the Kafka message-type test data.
# "JSONMessageType" -> "DemoEvent" (the verify-step's message type column)
# and introduce a new "JSONType" value for the identifier column that
# previously reused "DemoEvent".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API-KAKFA-PROTOBUFF-Testing")

$ws.Range("H2").Value = "DemoEvent"
$ws.Range("H6").Value = "DemoEvent"
$ws.Range("H7").Value = "DemoEvent"

$ws.Range("J2").Value = "JSONType"
$ws.Range("J6").Value = "JSONType"
$ws.Range("J7").Value = "JSONType"

# Move the active selection to J7, matching the saved view state.
$ws.Range("J7").Select()
